# Applies the "Added LR1, LR2, LR3" edit:
#  - Replaces the list of tenant surnames in column B (rows 3-38) with a new
#    roster (the old roster's shared strings become unused and are dropped;
#    new ones are appended to the shared string table).
#  - Changes the per-row tariff formula in column D from a hard-coded
#    constant (18.7 / 18.7/2) to a formula driven by $A$1 ($A$1*1.1 /
#    $A$1*0.55), which also changes the cached values of D, E, K, C40, C43.
#  - Updates the active window selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- Column B: new tenant surnames for rows 3-38 -----------------------
$ws.Range("B3").Value = "Абделазиз "
$ws.Range("B4").Value = ".Абдуллина "
$ws.Range("B5").Value = "Акмалов"
$ws.Range("B6").Value = "Бабкин "
$ws.Range("B7").Value = "Бахромов "
$ws.Range("B8").Value = "Борисов "
$ws.Range("B9").Value = "Гафеев "
$ws.Range("B10").Value = "Городилова "
$ws.Range("B11").Value = "Девятов   "
$ws.Range("B12").Value = "Исламов "
$ws.Range("B13").Value = "Карманов "
$ws.Range("B14").Value = "КОканов  "
$ws.Range("B15").Value = "МохамедБоуйе "
$ws.Range("B16").Value = "Надеждин  "
$ws.Range("B17").Value = "Нуретдинов  "
$ws.Range("B18").Value = "Павлова "
$ws.Range("B19").Value = "Родионов "
$ws.Range("B20").Value = "Рудой "
$ws.Range("B21").Value = "Садыков "
$ws.Range("B22").Value = "Семагин "
$ws.Range("B23").Value = "Семенов  "
$ws.Range("B24").Value = "Таухутдинов "
$ws.Range("B25").Value = "Фаляхутдинова "
$ws.Range("B26").Value = "Филатов  "
$ws.Range("B27").Value = "Хабк Осама "
$ws.Range("B28").Value = "Хазипова "
$ws.Range("B29").Value = "Хазов "
$ws.Range("B30").Value = "Хакимов "
$ws.Range("B31").Value = "Халилов "
$ws.Range("B32").Value = "Хафизов "
$ws.Range("B33").Value = "Хрунин "
$ws.Range("B34").Value = "Чепурченко  "
$ws.Range("B35").Value = "ШабАнов "
$ws.Range("B36").Value = "Шаймарданова "
$ws.Range("B37").Value = "Куропаткин 1"
$ws.Range("B38").Value = "Куропаткин 2"

# --- Column D: tariff formula now references $A$1 instead of a constant -
# D3 keeps a standalone formula, D4:D38 form one shared-formula group, just
# like in the original workbook.
$ws.Range("D3").Formula = "=IF(A3<=32,`$A`$1*1.1,`$A`$1*0.55)"
$ws.Range("D4:D38").FormulaR1C1 = "=IF(RC[-3]<=32,R1C1*1.1,R1C1*0.55)"

# --- Active window state -------------------------------------------------
$ws.Activate()
$ws.Range("D10").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
